$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1750
$ws.Range("J43").Value = 1750
$ws.Range("L43").Value = 1750
$ws.Range("N43").Value = -1888
$ws.Range("H62").Value = 7765.706
$ws.Range("J62").Value = 12812.375
$ws.Range("L62").Value = 12812.375
$ws.Range("N62").Value = -14060.375
$ws.Range("H65").Value = 7765.706
$ws.Range("J65").Value = 12812.375
$ws.Range("L65").Value = 64061.875
$ws.Range("N65").Value = -70301.875
$ws.Range("H96").Value = 420.83334
$ws.Range("I96").Value = 429
$ws.Range("K96").Value = 1287
$ws.Range("M96").Value = 86
$ws.Range("H99").Value = 1431.909
$ws.Range("I99").Value = 438
$ws.Range("K99").Value = 1314
$ws.Range("M99").Value = 184
$ws.Range("H101").Value = 443
$ws.Range("I101").Value = 419
$ws.Range("J101").Value = 467
$ws.Range("K101").Value = 1257
$ws.Range("L101").Value = 1401
$ws.Range("M101").Value = 365
$ws.Range("N101").Value = -4645
$ws.Range("H118").Value = 918.9
$ws.Range("J118").Value = 2333.3333
$ws.Range("L118").Value = 6999.999899999999
$ws.Range("N118").Value = -10313.9999
$ws.Range("H123").Value = 48999
$ws.Range("J123").Value = 48999
$ws.Range("L123").Value = 48999
$ws.Range("N123").Value = -58799
$ws.Range("H138").Value = 3005.6072
$ws.Range("J138").Value = 3586.7222
$ws.Range("L138").Value = 10760.1666
$ws.Range("N138").Value = -21040.1666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3200.05
$ws.Range("I2").Value = 3068.4614
$ws.Range("J2").Value = 3444.4285
$ws.Range("K2").Value = 3068.4614
$ws.Range("L2").Value = 3444.4285
$ws.Range("M2").Value = -2955.4614
$ws.Range("N2").Value = -3670.4285
$ws.Range("H61").Value = 1770.8235
$ws.Range("I61").Value = 1770.8235
$ws.Range("K61").Value = 1770.8235
$ws.Range("M61").Value = -1558.8235
$ws.Range("H116").Value = 3200.05
$ws.Range("I116").Value = 3068.4614
$ws.Range("J116").Value = 3444.4285
$ws.Range("K116").Value = 3068.4614
$ws.Range("L116").Value = 3444.4285
$ws.Range("M116").Value = -774.4614000000001
$ws.Range("N116").Value = -8032.4285
$ws.Range("H136").Value = 1770.8235
$ws.Range("I136").Value = 1770.8235
$ws.Range("K136").Value = 5312.470499999999
$ws.Range("M136").Value = -2762.470499999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3200.05
$ws.Range("I3").Value = 3068.4614
$ws.Range("J3").Value = 3444.4285
$ws.Range("K3").Value = 3068.4614
$ws.Range("L3").Value = 3444.4285
$ws.Range("M3").Value = -2954.4614
$ws.Range("N3").Value = -3672.4285
$ws.Range("H94").Value = 499.8889
$ws.Range("I94").Value = 400
$ws.Range("K94").Value = 400
$ws.Range("M94").Value = 51

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1270.8572
$ws.Range("I16").Value = 1074
$ws.Range("K16").Value = 1074
$ws.Range("M16").Value = -787
$ws.Range("H99").Value = 3562.4614
$ws.Range("I99").Value = 2888.7778
$ws.Range("J99").Value = 5078.25
$ws.Range("K99").Value = 2888.7778
$ws.Range("L99").Value = 5078.25
$ws.Range("M99").Value = -1390.7778
$ws.Range("N99").Value = -8074.25
$ws.Range("H103").Value = 20000
$ws.Range("I103").Value = 20000
$ws.Range("K103").Value = 20000
$ws.Range("M103").Value = -18828
$ws.Range("H113").Value = 1270.8572
$ws.Range("I113").Value = 1074
$ws.Range("K113").Value = 1074
$ws.Range("M113").Value = 1096
$ws.Range("H126").Value = 3562.4614
$ws.Range("I126").Value = 2888.7778
$ws.Range("J126").Value = 5078.25
$ws.Range("K126").Value = 8666.3334
$ws.Range("L126").Value = 15234.75
$ws.Range("M126").Value = -6196.3334
$ws.Range("N126").Value = -20174.75
$ws.Range("H134").Value = 2245.5557
$ws.Range("I134").Value = 2313.75
$ws.Range("J134").Value = 1700
$ws.Range("K134").Value = 6941.25
$ws.Range("L134").Value = 5100
$ws.Range("M134").Value = -4406.25
$ws.Range("N134").Value = -10170

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 949.6667
$ws.Range("I68").Value = 450
$ws.Range("J68").Value = 1199.5
$ws.Range("K68").Value = 1350
$ws.Range("L68").Value = 3598.5
$ws.Range("M68").Value = -539
$ws.Range("N68").Value = -5220.5
$ws.Range("H71").Value = 949.6667
$ws.Range("I71").Value = 450
$ws.Range("J71").Value = 1199.5
$ws.Range("K71").Value = 4050
$ws.Range("L71").Value = 10795.5
$ws.Range("M71").Value = 6
$ws.Range("N71").Value = -18907.5
$ws.Range("H131").Value = 3881.7
$ws.Range("J131").Value = 3910.875
$ws.Range("L131").Value = 11732.625
$ws.Range("N131").Value = -21812.625

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 49130.832
$ws.Range("J57").Value = 89595
$ws.Range("L57").Value = 89595
$ws.Range("N57").Value = -91235
$ws.Range("H113").Value = 5187.4
$ws.Range("I113").Value = 3127.7144
$ws.Range("J113").Value = 9993.333000000001
$ws.Range("K113").Value = 3127.7144
$ws.Range("L113").Value = 9993.333000000001
$ws.Range("M113").Value = -957.7143999999998
$ws.Range("N113").Value = -14333.333
$ws.Range("H122").Value = 2689.5334
$ws.Range("J122").Value = 5747.5
$ws.Range("L122").Value = 17242.5
$ws.Range("N122").Value = -22142.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3750
$ws.Range("J61").Value = 6999.6665
$ws.Range("L61").Value = 6999.6665
$ws.Range("N61").Value = -7403.6665
$ws.Range("H110").Value = 30000
$ws.Range("J110").Value = 30000
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
$ws.Range("H113").Value = 3750
$ws.Range("J113").Value = 6999.6665
$ws.Range("L113").Value = 6999.6665
$ws.Range("N113").Value = -11339.6665

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 31885.25
$ws.Range("I52").Value = 15021
$ws.Range("K52").Value = 15021
$ws.Range("M52").Value = -14795
$ws.Range("H113").Value = 1288.3334
$ws.Range("I113").Value = 1246
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 3738
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -1568
$ws.Range("N113").Value = -8840
